# First implementation with multiple AI providers
# The workbook cell A1 held the result of the initial (OpenAI) prompt run;
# it is updated here with the response obtained from the latest provider
# (Ollama / llama3.1) for the saved prompt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Ollama (llama3.1) response to 'Prova di un prompt SALVATO'"
